# Add new columns I (I0) and J (IF) to the worksheet, mirroring the
# existing header style (copied from column H) and filling in the
# per-row values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the bold/bordered/centered header style used by the other
# header cells (copy from H1, which already carries it).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Data rows 2-36: I = I0, J = IF values from the source data.
$values = @(
  @(2, 8, 8),
  @(3, 6, 6),
  @(4, 8, 8),
  @(5, 8, 8),
  @(6, 8, 8),
  @(7, 6, 6),
  @(8, 6, 6),
  @(9, 7, 9),
  @(10, 4, 5),
  @(11, 9, 9),
  @(12, 9, 9),
  @(13, 5, 5),
  @(14, 5, 5),
  @(15, 4, 5),
  @(16, 7, 8),
  @(17, 7, 7),
  @(18, 6, 6),
  @(19, 6, 7),
  @(20, 7, 7),
  @(21, 7, 8),
  @(22, 6, 6),
  @(23, 8, 8),
  @(24, 8, 8),
  @(25, 9, 9),
  @(26, 9, 9),
  @(27, 2, 3),
  @(28, 8, 8),
  @(29, 1, 2),
  @(30, 4, 4),
  @(31, 8, 9),
  @(32, 4, 5),
  @(33, 9, 9),
  @(34, 8, 8),
  @(35, 4, 4),
  @(36, 5, 5)
)

foreach ($row in $values) {
  $r = $row[0]
  $ws.Cells.Item($r, 9).Value = $row[1]
  $ws.Cells.Item($r, 10).Value = $row[2]
}
